$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the course title for A5 (was "AFAM 102.001", should be "AFAM 102C.001")
$ws.Range("A5").Value = "AFAM 102C.001"

# Update the active selection to A5 (matches the saved sheet view state)
$ws.Range("A5").Select() | Out-Null
